# Update vm_pu results for the 380 kV case (Case_2_100)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.017347447144331
$rowBF[0,2] = 1.023252758806419
$rowBF[0,3] = 1.04485218078264
$rowBF[0,4] = 1.048072276014635
$ws.Range("B2:F2").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.026900141761679
$rowIN[0,1] = 1.022562078628814
$rowIN[0,2] = 1.026084741980745
$rowIN[0,3] = 1.047622066623595
$rowIN[0,4] = 1.050833142073923
$rowIN[0,5] = 1.011623335092795
$ws.Range("I2:N2").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.018192036312429
$rowBF[0,2] = 1.023856007039742
$rowBF[0,3] = 1.046083947833962
$rowBF[0,4] = 1.049365127179596
$ws.Range("B3:F3").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.026985701210567
$rowIN[0,1] = 1.023043066375941
$rowIN[0,2] = 1.026495385174521
$rowIN[0,3] = 1.048664059751691
$rowIN[0,4] = 1.051936720123
$rowIN[0,5] = 1.011783283761708
$ws.Range("I3:N3").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.018738763131138
$rowBF[0,2] = 1.024246269601081
$rowBF[0,3] = 1.046881853694204
$rowBF[0,4] = 1.050202610550996
$ws.Range("B4:F4").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.02703963978873
$rowIN[0,1] = 1.023353894167329
$rowIN[0,2] = 1.026760326248456
$rowIN[0,3] = 1.049338607238758
$rowIN[0,4] = 1.052651181048792
$rowIN[0,5] = 1.011886614886638
$ws.Range("I4:N4").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.018968659008877
$rowBF[0,2] = 1.024410315520783
$rowBF[0,3] = 1.047217502670192
$rowBF[0,4] = 1.050554910351702
$ws.Range("B5:F5").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.02706197409507
$rowIN[0,1] = 1.023484468789004
$rowIN[0,2] = 1.026871521537065
$rowIN[0,3] = 1.049622261644932
$rowIN[0,4] = 1.052951630282977
$rowIN[0,5] = 1.011930015161328
$ws.Range("I5:N5").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.019007262553479
$rowBF[0,2] = 1.024437858310245
$rowBF[0,3] = 1.047273871916593
$rowBF[0,4] = 1.050614076073954
$ws.Range("B6:F6").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.027065704083649
$rowIN[0,1] = 1.023506387108838
$rowIN[0,2] = 1.026890180782827
$rowIN[0,3] = 1.049669892869678
$rowIN[0,4] = 1.053002082323512
$rowIN[0,5] = 1.011937299896323
$ws.Range("I6:N6").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.018741834808365
$rowBF[0,2] = 1.024248461672929
$rowBF[0,3] = 1.046886337828345
$rowBF[0,4] = 1.050207317125408
$ws.Range("B7:F7").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.027039939563251
$rowIN[0,1] = 1.023355639295463
$rowIN[0,2] = 1.026761812777213
$rowIN[0,3] = 1.049342397148454
$rowIN[0,4] = 1.052655195312098
$rowIN[0,5] = 1.0118871949613
$ws.Range("I7:N7").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.017632833603589
$rowBF[0,2] = 1.023456644934078
$rowBF[0,3] = 1.045268282654309
$rowBF[0,4] = 1.0485090107738
$ws.Range("B8:F8").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.026929351525518
$rowIN[0,1] = 1.02272471374725
$rowIN[0,2] = 1.026223680132148
$rowIN[0,3] = 1.047974149759866
$rowIN[0,4] = 1.051206024852075
$rowIN[0,5] = 1.011677424757957
$ws.Range("I8:N8").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.015680370445265
$rowBF[0,2] = 1.022060818779864
$rowBF[0,3] = 1.042423681797867
$rowBF[0,4] = 1.045523395857631
$ws.Range("B9:F9").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.026723598010628
$rowIN[0,1] = 1.021609887423023
$rowIN[0,2] = 1.025269548670582
$rowIN[0,3] = 1.045565446112952
$rowIN[0,4] = 1.04865521688417
$rowIN[0,5] = 1.011306522076245
$ws.Range("I9:N9").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.014379956500186
$rowBF[0,2] = 1.021129988111977
$rowBF[0,3] = 1.040531665623966
$rowBF[0,4] = 1.043537638909543
$ws.Range("B10:F10").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.026579143363579
$rowIN[0,1] = 1.020864662243278
$rowIN[0,2] = 1.024629570583584
$rowIN[0,3] = 1.043961150843651
$rowIN[0,4] = 1.046956520806797
$rowIN[0,5] = 1.011058423630407
$ws.Range("I10:N10").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.013817166248488
$rowBF[0,2] = 1.020726879665883
$rowBF[0,3] = 1.039713426337722
$rowBF[0,4] = 1.042678872596296
$ws.Range("B11:F11").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.026514871417588
$rowIN[0,1] = 1.020541504808619
$rowIN[0,2] = 1.02435154229953
$rowIN[0,3] = 1.043266817974231
$rowIN[0,4] = 1.046221391777876
$rowIN[0,5] = 1.010950801037839
$ws.Range("I11:N11").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.013608166513698
$rowBF[0,2] = 1.020577140755195
$rowBF[0,3] = 1.039409647106958
$rowBF[0,4] = 1.042360049339323
$ws.Range("B12:F12").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.026490739725249
$rowIN[0,1] = 1.020421399859451
$rowIN[0,2] = 1.024248133902741
$rowIN[0,3] = 1.043008961365806
$rowIN[0,4] = 1.045948393748533
$rowIN[0,5] = 1.010910796299761
$ws.Range("I12:N12").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.013652995587025
$rowBF[0,2] = 1.020609260530565
$rowBF[0,3] = 1.039474801960205
$rowBF[0,4] = 1.042428430756653
$ws.Range("B13:F13").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.026495927732588
$rowIN[0,1] = 1.020447165929774
$rowIN[0,2] = 1.024270321508306
$rowIN[0,3] = 1.043064270261004
$rowIN[0,4] = 1.046006949981765
$rowIN[0,5] = 1.010919378755866
$ws.Range("I13:N13").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.013799889329216
$rowBF[0,2] = 1.020714502324919
$rowBF[0,3] = 1.039688312773915
$rowBF[0,4] = 1.042652515294469
$ws.Range("B14:F14").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.02651288194929
$rowIN[0,1] = 1.02053157831759
$rowIN[0,2] = 1.024342997304524
$rowIN[0,3] = 1.043245502462865
$rowIN[0,4] = 1.046198824424925
$rowIN[0,5] = 1.010947494821956
$ws.Range("I14:N14").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.013890401484792
$rowBF[0,2] = 1.020779344450856
$rowBF[0,3] = 1.039819883836235
$rowBF[0,4] = 1.042790602438132
$ws.Range("B15:F15").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.02652329379876
$rowIN[0,1] = 1.02058357833632
$rowIN[0,2] = 1.024387757226386
$rowIN[0,3] = 1.043357172132404
$rowIN[0,4] = 1.046317052734532
$rowIN[0,5] = 1.010964814235031
$ws.Range("I15:N15").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.014417313383344
$rowBF[0,2] = 1.021156740104617
$rowBF[0,3] = 1.040585990761116
$rowBF[0,4] = 1.043594655026555
$ws.Range("B16:F16").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.026583372640409
$rowIN[0,1] = 1.020886099297188
$rowIN[0,2] = 1.024648003237803
$rowIN[0,3] = 1.04400723844068
$rowIN[0,4] = 1.047005317602176
$rowIN[0,5] = 1.011065562117107
$ws.Range("I16:N16").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.014747911611043
$rowBF[0,2] = 1.021393457343414
$rowBF[0,3] = 1.041066820470889
$rowBF[0,4] = 1.044099304007535
$ws.Range("B17:F17").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.026620597801025
$rowIN[0,1] = 1.021075737379735
$rowIN[0,2] = 1.024811004765018
$rowIN[0,3] = 1.044415097402705
$rowIN[0,4] = 1.047437158981242
$rowIN[0,5] = 1.011128706772829
$ws.Range("I17:N17").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.014940772613499
$rowBF[0,2] = 1.021531525331413
$rowBF[0,3] = 1.041347378448042
$rowBF[0,4] = 1.044393761474703
$ws.Range("B18:F18").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.026642144413243
$rowIN[0,1] = 1.021186304666802
$rowIN[0,2] = 1.024905992591122
$rowIN[0,3] = 1.044653027284265
$rowIN[0,4] = 1.047689085116479
$rowIN[0,5] = 1.011165519213579
$ws.Range("I18:N18").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.015006538062927
$rowBF[0,2] = 1.021578602026675
$rowBF[0,3] = 1.041443058157334
$rowBF[0,4] = 1.04449418158765
$ws.Range("B19:F19").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.026649463046535
$rowIN[0,1] = 1.021223997531996
$rowIN[0,2] = 1.024938365995119
$rowIN[0,3] = 1.044734160799152
$rowIN[0,4] = 1.047774992298796
$rowIN[0,5] = 1.01117806811253
$ws.Range("I19:N19").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.014712438579934
$rowBF[0,2] = 1.02136806034868
$rowBF[0,3] = 1.041015221852147
$rowBF[0,4] = 1.044045149174211
$ws.Range("B20:F20").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.026616621080201
$rowIN[0,1] = 1.021055395698655
$rowIN[0,2] = 1.024793525360635
$rowIN[0,3] = 1.044371334627296
$rowIN[0,4] = 1.047390822320898
$rowIN[0,5] = 1.011121933888152
$ws.Range("I20:N20").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.013756631516591
$rowBF[0,2] = 1.020683511407211
$rowBF[0,3] = 1.039625434996259
$rowBF[0,4] = 1.042586523555925
$ws.Range("B21:F21").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.026507896479435
$rowIN[0,1] = 1.020506722896776
$rowIN[0,2] = 1.024321599849322
$rowIN[0,3] = 1.043192132738578
$rowIN[0,4] = 1.04614232049156
$rowIN[0,5] = 1.010939216137224
$ws.Range("I21:N21").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.013155942094173
$rowBF[0,2] = 1.020253071494576
$rowBF[0,3] = 1.038752495364843
$rowBF[0,4] = 1.041670357319081
$ws.Range("B22:F22").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.026438043216497
$rowIN[0,1] = 1.020161346528809
$rowIN[0,2] = 1.024024093347198
$rowIN[0,3] = 1.042451007718959
$rowIN[0,4] = 1.045357693826575
$rowIN[0,5] = 1.010824167078206
$ws.Range("I22:N22").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.013474353552867
$rowBF[0,2] = 1.02048125885402
$rowBF[0,3] = 1.039215174751476
$rowBF[0,4] = 1.04215594672872
$ws.Range("B23:F23").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.026475215179987
$rowIN[0,1] = 1.020344475129224
$rowIN[0,2] = 1.024181881573269
$rowIN[0,3] = 1.042843865397133
$rowIN[0,4] = 1.045773605894255
$rowIN[0,5] = 1.010885172512399
$ws.Range("I23:N23").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.0147284672273
$rowBF[0,2] = 1.021379536173526
$rowBF[0,3] = 1.04103853674327
$rowBF[0,4] = 1.044069619089777
$ws.Range("B24:F24").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.026618418502722
$rowIN[0,1] = 1.021064587369054
$rowIN[0,2] = 1.024801423823706
$rowIN[0,3] = 1.044391109040518
$rowIN[0,4] = 1.047411759738094
$rowIN[0,5] = 1.011124994321214
$ws.Range("I24:N24").Value = $rowIN

$rowBF = New-Object "object[,]" 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.016184916833098
$rowBF[0,2] = 1.022421728892203
$rowBF[0,3] = 1.043158301717278
$rowBF[0,4] = 1.046294423985439
$ws.Range("B25:F25").Value = $rowBF

$rowIN = New-Object "object[,]" 1,6
$rowIN[0,0] = 1.026778076684119
$rowIN[0,1] = 1.021898453992095
$rowIN[0,2] = 1.025516904492058
$rowIN[0,3] = 1.046187884104652
$rowIN[0,4] = 1.049314333445733
$rowIN[0,5] = 1.011402556892834
$ws.Range("I25:N25").Value = $rowIN

